{"js": "// Replace the date line and each two-digit-by-two-digit multiplication\n// expression in the table with its updated value, as described by the diff.\n// Each (oldText -> newText) pair is unique within the document, so a simple\n// exact-text search-and-replace for each pair is unambiguous.\nconst replacements = [\n  [\"2025-01-04 Saturday\", \"2025-01-05 Sunday\"],\n  [\"37\u00d756=\", \"37\u00d742=\"],\n  [\"13\u00d776=\", \"85\u00d714=\"],\n  [\"96\u00d738=\", \"96\u00d748=\"],\n  [\"11\u00d731=\", \"64\u00d740=\"],\n  [\"61\u00d781=\", \"77\u00d743=\"],\n  [\"97\u00d757=\", \"70\u00d782=\"],\n  [\"25\u00d729=\", \"40\u00d759=\"],\n  [\"16\u00d712=\", \"67\u00d745=\"],\n  [\"31\u00d764=\", \"45\u00d751=\"],\n  [\"34\u00d789=\", \"78\u00d717=\"],\n  [\"89\u00d716=\", \"42\u00d746=\"],\n  [\"12\u00d715=\", \"83\u00d739=\"],\n  [\"86\u00d788=\", \"48\u00d752=\"],\n  [\"81\u00d795=\", \"86\u00d796=\"],\n  [\"67\u00d773=\", \"39\u00d712=\"],\n  [\"61\u00d747=\", \"49\u00d721=\"],\n  [\"59\u00d714=\", \"18\u00d731=\"],\n  [\"59\u00d784=\", \"98\u00d794=\"],\n  [\"55\u00d749=\", \"21\u00d716=\"],\n  [\"55\u00d783=\", \"27\u00d794=\"],\n  [\"68\u00d719=\", \"87\u00d752=\"],\n  [\"79\u00d736=\", \"11\u00d775=\"],\n  [\"31\u00d728=\", \"20\u00d739=\"],\n  [\"24\u00d773=\", \"96\u00d759=\"],\n  [\"66\u00d749=\", \"12\u00d722=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit-by-two-digit multiplication\n# expression in the table with its updated value, as described by the diff.\n# Each old value is unique in the document, so Find/Replace (wdReplaceAll)\n# for each pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-01-04 Saturday\", \"2025-01-05 Sunday\"),\n    @(\"37\u00d756=\", \"37\u00d742=\"),\n    @(\"13\u00d776=\", \"85\u00d714=\"),\n    @(\"96\u00d738=\", \"96\u00d748=\"),\n    @(\"11\u00d731=\", \"64\u00d740=\"),\n    @(\"61\u00d781=\", \"77\u00d743=\"),\n    @(\"97\u00d757=\", \"70\u00d782=\"),\n    @(\"25\u00d729=\", \"40\u00d759=\"),\n    @(\"16\u00d712=\", \"67\u00d745=\"),\n    @(\"31\u00d764=\", \"45\u00d751=\"),\n    @(\"34\u00d789=\", \"78\u00d717=\"),\n    @(\"89\u00d716=\", \"42\u00d746=\"),\n    @(\"12\u00d715=\", \"83\u00d739=\"),\n    @(\"86\u00d788=\", \"48\u00d752=\"),\n    @(\"81\u00d795=\", \"86\u00d796=\"),\n    @(\"67\u00d773=\", \"39\u00d712=\"),\n    @(\"61\u00d747=\", \"49\u00d721=\"),\n    @(\"59\u00d714=\", \"18\u00d731=\"),\n    @(\"59\u00d784=\", \"98\u00d794=\"),\n    @(\"55\u00d749=\", \"21\u00d716=\"),\n    @(\"55\u00d783=\", \"27\u00d794=\"),\n    @(\"68\u00d719=\", \"87\u00d752=\"),\n    @(\"79\u00d736=\", \"11\u00d775=\"),\n    @(\"31\u00d728=\", \"20\u00d739=\"),\n    @(\"24\u00d773=\", \"96\u00d759=\"),\n    @(\"66\u00d749=\", \"12\u00d722=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $oldText,    # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw \"No match found for: \" + $oldText\n    }\n}\n\n"}
